$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

Set-TextValue "D2" '58.050.01'
Set-TextValue "E2" '  -1.37%  '
Set-TextValue "D3" '2.452.08'
Set-TextValue "E3" '  -3.67%  '
Set-TextValue "E4" '  +0.24%  '
Set-TextValue "D5" '524.63'
Set-TextValue "E5" '  -0.47%  '
Set-TextValue "D6" '130.44'
Set-TextValue "E6" '  -2.69%  '
Set-TextValue "E7" '  +0.18%  '
Set-TextValue "E8" '  -0.28%  '
Set-TextValue "D9" '2.454.00'
Set-TextValue "E9" '  -3.38%  '
Set-TextValue "E10" '  -0.81%  '
Set-TextValue "E11" '  -2.12%  '
Set-TextValue "D12" '4.97'
Set-TextValue "E12" '  -4.16%  '
Set-TextValue "E13" '  -2.92%  '
Set-TextValue "D14" '2.890.86'
Set-TextValue "E14" '  -3.46%  '
Set-TextValue "D15" '57.967.27'
Set-TextValue "E15" '  -1.36%  '
Set-TextValue "D16" '21.66'
Set-TextValue "E16" '  -3.04%  '
Set-TextValue "E17" '  -2.23%  '
Set-TextValue "D18" '2.452.47'
Set-TextValue "E18" '  -3.36%  '
Set-TextValue "E19" '  -3.09%  '
Set-TextValue "E20" '  -1.67%  '
Set-TextValue "D21" '314.63'
Set-TextValue "E21" '  -2.99%  '
Set-TextValue "D22" '6.14'
Set-TextValue "E22" '  -0.62%  '
Set-TextValue "E23" '  -0.02%  '
Set-TextValue "D24" '65.23'
Set-TextValue "E24" '  -0.26%  '
Set-TextValue "E25" '  -1.35%  '
Set-TextValue "D26" '2.568.18'
Set-TextValue "E26" '  -2.84%  '
Set-TextValue "E27" '  +0.71%  '
Set-TextValue "E28" '  -1.93%  '
Set-TextValue "E29" '  -2.22%  '
Set-TextValue "D30" '173.52'
Set-TextValue "E30" '  +3.03%  '
Set-TextValue "E31" '  -2.75%  '
Set-TextValue "D32" '1.70'
Set-TextValue "E32" '  -2.36%  '
Set-TextValue "E33" '  -3.02%  '
Set-TextValue "E34" '  -6.12%  '
Set-TextValue "D35" '0.999'
Set-TextValue "E35" '  -0.01%  '
Set-TextValue "D36" '0.996'
Set-TextValue "E36" '  -0.20%  '
Set-TextValue "E37" '  -2.41%  '
Set-TextValue "E38" '  -6.58%  '
Set-TextValue "E39" '  -4.33%  '
Set-TextValue "E40" '  -0.63%  '
Set-TextValue "D41" '0.813'
Set-TextValue "E41" '  +3.48%  '
Set-TextValue "E42" '  -3.00%  '
Set-TextValue "E43" '  -2.32%  '
Set-TextValue "B44" 'Mantle'
Set-TextValue "C44" 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue "D44" '0.584'
Set-TextValue "E44" '  -3.49%  '
Set-TextValue "B45" 'Bittensor'
Set-TextValue "C45" 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue "D45" '260.22'
Set-TextValue "E45" '  -7.20%  '
Set-TextValue "D46" '4.81'
Set-TextValue "E46" '  -7.07%  '
Set-TextValue "B47" 'Aave'
Set-TextValue "C47" 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue "D47" '124.18'
Set-TextValue "E47" '  -4.66%  '
Set-TextValue "B48" 'Stellar'
Set-TextValue "C48" 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue "D48" '0.0922'
Set-TextValue "E48" '  +0.34%  '
Set-TextValue "E49" '  -2.24%  '
Set-TextValue "E50" '  -2.14%  '
Set-TextValue "D51" '17.04'
Set-TextValue "E51" '  -4.86%  '
